$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14: add a new progress-report entry (date, hours, activity),
# matching the date-cell formatting already used by rows 3-13.
$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A14").Value = 42796         # 3/2/2017

$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "Tested and Fixed SQL files"

# Update the active selection to B15 (matches the saved sheet view state)
$ws.Range("B15").Select()
